$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data currently sits on row 1 (header) and row 3 (data), with an
# empty row 2 in between. Remove the extra empty row so the data row
# shifts up to row 2.
$ws.Rows.Item(2).Delete()

# Update the selection to match the target state.
$ws.Range("E6").Select()
